# Katie - moved over stub fbx loader, updated milestone 1, added fbx objects
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename the FBX rubric rows (row 16 & 17) ---
# Row 16: "FBX mesh data to binary" / "Extract FBX mesh data into a binary file"
#     ->  "FBX mesh data "        / "Extract FBX mesh data"
$ws.Range("A16").Value = "FBX mesh data "
$ws.Range("C16").Value = "Extract FBX mesh data"

# Row 17: "FBX render Binary" / "Render meshes from Binary files"
#     ->  "FBX render"        / "Render meshes "
$ws.Range("A17").Value = "FBX render"
$ws.Range("C17").Value = "Render meshes "

# --- B8 should be bold + right aligned like the other point-value cells (B10/B11/B18) ---
$ws.Range("B8").Font.Bold = $true

# --- Update the sheet view: scroll down a bit and move the selection to C25 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C25").Select()
